$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final sorted data (label, value) descending by value, with the new
# "bogota d.c. (ct)" row inserted and "cundinamarca" updated/re-sorted.
$data = @(
  @("bogota d.c. (ct)", 54323),
  @("antioquia", 39636),
  @("valle", 32245),
  @("santander", 20504),
  @("cundinamarca", 16338),
  @("atlantico", 14269),
  @("huila", 12485),
  @("tolima", 11902),
  @("bolivar", 11341),
  @("meta", 10317),
  @("boyaca", 9874),
  @("narino", 9250),
  @("risaralda", 9220),
  @("cauca", 8601),
  @("cordoba", 7724),
  @("norte de santander", 7628),
  @("caldas", 7626),
  @("magdalena", 7442),
  @("cesar", 6922),
  @("sucre", 5528),
  @("caqueta", 4308),
  @("quindio", 4215),
  @("guajira", 3493),
  @("putumayo", 3114),
  @("casanare", 3030),
  @("choco", 2554),
  @("arauca", 2441),
  @("guaviare", 1174),
  @("amazonas", 1080),
  @("san andres", 631),
  @("vichada", 489),
  @("guainia", 407),
  @("vaupes", 276)
)

$row = 2
foreach ($pair in $data) {
  $ws.Cells.Item($row, 1).Value = $pair[0]
  $ws.Cells.Item($row, 2).Value = $pair[1]
  $row = $row + 1
}
$lastRow = $row - 1

# Apply left-horizontal-alignment style to all the department cells (A2:A34)
$ws.Range("A2:A$lastRow").HorizontalAlignment = -4131

# Update the active selection shown when the workbook is opened
$ws.Range("C11").Select()
